# Feat: Show MaxAttackRange Decal
# Insert a new "MaxBasicAttackRange" column before the existing "BaseDamage" column (E),
# shifting BaseDamage..HpRegenerationRate one column to the right (F..L).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GOCharacterStatTable")

# Insert a new column at E; existing E:K shift to F:L.
$ws.Range("E1").EntireColumn.Insert()

# Header for the new column.
$ws.Range("E1").Value = "MaxBasicAttackRange"

# Values for the new column (MaxBasicAttackRange) per character row.
$ws.Range("E2").Value = 75
$ws.Range("E3").Value = 350
$ws.Range("E4").Value = 150
$ws.Range("E5").Value = 250

# Width for the newly inserted column (the other columns keep their widths
# automatically, since EntireColumn.Insert() shifts the existing per-column
# widths along with the columns). 21.0 is the nearest value this engine's
# character->pixel quantization (MDW 7, 5px padding) can resolve to the
# target stored width of 21.6640625.
$ws.Columns.Item(5).ColumnWidth = 21
